$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab to reflect new "through" date
$ws.Name = "Through 2022-08-15"

# Update row label for August to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-15)"

# Update July 2022 (I7) value
$ws.Cells.Item(7, 9).Value = 143

# Update August row (row 9) values for years 2015-2022 (columns B-I)
$ws.Cells.Item(9, 2).Value = 15
$ws.Cells.Item(9, 3).Value = 35
$ws.Cells.Item(9, 4).Value = 34
$ws.Cells.Item(9, 5).Value = 24
$ws.Cells.Item(9, 6).Value = 20
$ws.Cells.Item(9, 7).Value = 93
$ws.Cells.Item(9, 8).Value = 87
$ws.Cells.Item(9, 9).Value = 85

# Update Total row (row 10) values for years 2015-2022 (columns B-I)
$ws.Cells.Item(10, 2).Value = 177
$ws.Cells.Item(10, 3).Value = 337
$ws.Cells.Item(10, 4).Value = 499
$ws.Cells.Item(10, 5).Value = 449
$ws.Cells.Item(10, 6).Value = 324
$ws.Cells.Item(10, 7).Value = 714
$ws.Cells.Item(10, 8).Value = 997
$ws.Cells.Item(10, 9).Value = 1056
